$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.876.44"
$ws.Range("E2").Value = "  -5.30%  "
$ws.Range("D3").Value = "3.217.13"
$ws.Range("E3").Value = "  -6.53%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'174.13"
$ws.Range("E5").Value = "  -7.21%  "
$ws.Range("D6").Value = "'513.63"
$ws.Range("E6").Value = "  -5.10%  "
$ws.Range("D7").Value = "'0.590"
$ws.Range("E7").Value = "  -5.19%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "3.213.74"
$ws.Range("E9").Value = "  -6.51%  "
$ws.Range("E10").Value = "  -7.14%  "
$ws.Range("D11").Value = "'52.22"
$ws.Range("E11").Value = "  -12.15%  "
$ws.Range("E12").Value = "  -6.68%  "
$ws.Range("E13").Value = "  -4.16%  "
$ws.Range("E14").Value = "  -7.55%  "
$ws.Range("D15").Value = "3.731.95"
$ws.Range("E15").Value = "  -6.09%  "
$ws.Range("E16").Value = "  -7.25%  "
$ws.Range("D17").Value = "3.213.56"
$ws.Range("E17").Value = "  -6.54%  "
$ws.Range("D18").Value = "62.838.50"
$ws.Range("E18").Value = "  -4.90%  "
$ws.Range("E19").Value = "  -4.36%  "
$ws.Range("D20").Value = "'10.89"
$ws.Range("E20").Value = "  -5.79%  "
$ws.Range("D21").Value = "'0.952"
$ws.Range("E21").Value = "  -4.99%  "
$ws.Range("D22").Value = "'363.68"
$ws.Range("E22").Value = "  -6.30%  "
$ws.Range("E23").Value = "  -4.12%  "
$ws.Range("D24").Value = "'79.84"
$ws.Range("E24").Value = "  -5.19%  "
$ws.Range("D25").Value = "'10.96"
$ws.Range("E25").Value = "  -2.13%  "
$ws.Range("D26").Value = "'3.89"
$ws.Range("E26").Value = "  +1.63%  "
$ws.Range("D27").Value = "'6.04"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("E28").Value = "  -6.07%  "
$ws.Range("E29").Value = "  -7.13%  "
$ws.Range("D30").Value = "'8.11"
$ws.Range("E30").Value = "  -7.17%  "
$ws.Range("D31").Value = "'649.93"
$ws.Range("E31").Value = "  -7.08%  "
$ws.Range("E32").Value = "  -7.31%  "
$ws.Range("D33").Value = "'6.23"
$ws.Range("E33").Value = "  -10.23%  "
$ws.Range("E34").Value = "  -3.53%  "
$ws.Range("E35").Value = "  -4.79%  "
$ws.Range("D36").Value = "'57.66"
$ws.Range("E36").Value = "  -7.31%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "'36.27"
$ws.Range("E38").Value = "  -3.15%  "
$ws.Range("E39").Value = "  -5.25%  "
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").Value = "0.0₃0691"
$ws.Range("E41").Value = "  +7.78%  "
$ws.Range("E42").Value = "  -5.43%  "
$ws.Range("D43").Value = "2.850.44"
$ws.Range("E43").Value = "  -3.10%  "
$ws.Range("D44").Value = "'2.51"
$ws.Range("E44").Value = "  +3.21%  "
$ws.Range("D45").Value = "'2.67"
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").Value = "'2.85"
$ws.Range("E46").Value = "  +7.74%  "
$ws.Range("E47").Value = "  -3.07%  "
$ws.Range("D48").Value = "'2.57"
$ws.Range("E48").Value = "  -10.10%  "
$ws.Range("D49").Value = "'134.76"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").Value = "'2.90"
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.122"
$ws.Range("E51").Value = "  -4.70%  "
